# Add the new recipient e-mail address in the next empty row of column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "lyudvognova@gmail.com"

# Reflect the active selection recorded in the sheet view.
[void]$ws.Range("K12").Select()
